# Updated cryptos list refresh (prices / 1h volume changes), matching the
# "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price values (column D) are plain decimal numbers
# (e.g. "0.650", "1.00", "54.73"). If assigned as-is, Excel would
# auto-convert them to numeric values and normalize/trim them
# (e.g. "0.650" -> 0.65, "1.00" -> 1), which would NOT match the
# original text content of these cells. Forcing a text ("@") number
# format on those specific cells before writing keeps the values as
# literal text, exactly as stored in the source sheet.
$textCells = @(
    "D5", "D6", "D7", "D10", "D11", "D12", "D13", "D14", "D17", "D18",
    "D20", "D25", "D26", "D27", "D29", "D30", "D32", "D34", "D35", "D36",
    "D37", "D38", "D39", "D45", "D46", "D48", "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "70.727.26"
$ws.Range("E2").Value = "  +2.56%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "3.561.71"
$ws.Range("E3").Value = "  +1.46%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.05%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "581.54"
$ws.Range("E5").Value = "  +1.92%  "

# --- Row 6: Solana ---
$ws.Range("D6").Value = "187.96"
$ws.Range("E6").Value = "  +1.63%  "

# --- Row 7: XRP ---
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +2.30%  "

# --- Row 8: LidoStakedEther ---
$ws.Range("D8").Value = "3.551.37"
$ws.Range("E8").Value = "  +1.38%  "

# --- Row 9: USDC ---
$ws.Range("E9").Value = "  -0.08%  "

# --- Row 10: Dogecoin ---
$ws.Range("D10").Value = "0.223"
$ws.Range("E10").Value = "  +19.79%  "

# --- Row 11: Cardano ---
$ws.Range("D11").Value = "0.650"
$ws.Range("E11").Value = "  +0.02%  "

# --- Row 12: Avalanche ---
$ws.Range("D12").Value = "54.73"
$ws.Range("E12").Value = "  +1.03%  "

# --- Row 13: ShibaInu ---
$ws.Range("D13").Value = "0.0000320"
$ws.Range("E13").Value = "  +6.22%  "

# --- Row 14: Polkadot ---
$ws.Range("D14").Value = "9.52"
$ws.Range("E14").Value = "  +0.88%  "

# --- Row 15: WrappedliquidstakedEther2.0 ---
$ws.Range("D15").Value = "4.131.25"
$ws.Range("E15").Value = "  +1.43%  "

# --- Row 16: WrappedBTC ---
$ws.Range("D16").Value = "70.769.75"
$ws.Range("E16").Value = "  +2.68%  "

# --- Row 17: Chainlink ---
$ws.Range("D17").Value = "19.14"
$ws.Range("E17").Value = "  -1.08%  "

# --- Row 18: Uniswap ---
$ws.Range("D18").Value = "12.78"
$ws.Range("E18").Value = "  +4.34%  "

# --- Row 19: WrappedEther ---
$ws.Range("D19").Value = "3.558.52"
$ws.Range("E19").Value = "  +1.66%  "

# --- Row 20: BitcoinCash ---
$ws.Range("D20").Value = "576.61"
$ws.Range("E20").Value = "  +6.36%  "

# --- Row 22: Polygon ---
$ws.Range("E22").Value = "  -0.69%  "

# --- Row 23: InternetComputer(DFINITY) ---
$ws.Range("E23").Value = "  -4.54%  "

# --- Row 24: PancakeSwap ---
$ws.Range("E24").Value = "  +3.61%  "

# --- Row 25: Toncoin ---
$ws.Range("D25").Value = "4.86"
$ws.Range("E25").Value = "  -2.59%  "

# --- Row 26: Litecoin ---
$ws.Range("D26").Value = "93.82"
$ws.Range("E26").Value = "  -0.15%  "

# --- Row 27: RenderToken ---
$ws.Range("D27").Value = "11.24"
$ws.Range("E27").Value = "  +4.04%  "

# --- Row 28: ImmutableX ---
$ws.Range("E28").Value = "  +1.85%  "

# --- Row 29: Filecoin ---
$ws.Range("D29").Value = "9.29"
$ws.Range("E29").Value = "  +1.83%  "

# --- Row 30: EthereumClassic ---
$ws.Range("D30").Value = "32.57"
$ws.Range("E30").Value = "  +2.13%  "

# --- Row 31: NEARProtocol ---
$ws.Range("E31").Value = "  -0.55%  "

# --- Row 32: Cosmos ---
$ws.Range("D32").Value = "12.32"
$ws.Range("E32").Value = "  -1.88%  "

# --- Row 33: Hedera ---
$ws.Range("E33").Value = "  +1.97%  "

# --- Row 34 & 35: OKB / dogwifhat swap places in the ranking ---
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "63.11"
$ws.Range("E34").Value = "  -2.58%  "

$ws.Range("B35").Value = "dogwifhat"
$ws.Range("C35").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D35").Value = "3.74"
$ws.Range("E35").Value = "  +18.58%  "

# --- Row 36: Fetch.AI ---
$ws.Range("D36").Value = "3.33"
$ws.Range("E36").Value = "  +10.80%  "

# --- Row 37: Bittensor ---
$ws.Range("D37").Value = "541.37"
$ws.Range("E37").Value = "  -3.47%  "

# --- Row 38: TheGraph ---
$ws.Range("D38").Value = "0.411"
$ws.Range("E38").Value = "  +3.57%  "

# --- Row 39: InjectiveProtocol ---
$ws.Range("D39").Value = "38.25"
$ws.Range("E39").Value = "  +0.89%  "

# --- Row 40: PEPE ---
$ws.Range("D40").Value = "0.0₃0808"
$ws.Range("E40").Value = "  +5.68%  "

# --- Row 41: Dai ---
$ws.Range("E41").Value = "  -0.10%  "

# --- Row 42: Maker ---
$ws.Range("D42").Value = "3.571.91"
$ws.Range("E42").Value = "  +10.58%  "

# --- Row 43: Kaspa ---
$ws.Range("E43").Value = "  +4.43%  "

# --- Row 44: Stacks ---
$ws.Range("E44").Value = "  +2.93%  "

# --- Row 45: VeChain ---
$ws.Range("D45").Value = "0.0472"
$ws.Range("E45").Value = "  +7.40%  "

# --- Row 46: ApeXProtocol ---
$ws.Range("D46").Value = "3.50"
$ws.Range("E46").Value = "  -1.33%  "

# --- Row 47: ThetaToken ---
$ws.Range("E47").Value = "  -1.11%  "

# --- Row 48: THORChain ---
$ws.Range("D48").Value = "9.34"
$ws.Range("E48").Value = "  +4.34%  "

# --- Row 49: Stellar ---
$ws.Range("E49").Value = "  +2.83%  "

# --- Row 50: FirstDigitalUSD ---
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.15%  "

# --- Row 51: OceanProtocol ---
$ws.Range("E51").Value = "  +7.30%  "

# Restore the default (Normal) cell style on the cells where we forced a
# text number format, so the only persisted difference is the cell value
# itself - no stray style index is left referenced on these cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
